$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.102275
$ws.Range("H2").Value = 0.306825
$ws.Range("I2").Value = 0.2304482333616488
$ws.Range("J2").Value = 0.2304482333616488
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.594753
$ws.Range("N2").Value = 1.784259
$ws.Range("O2").Value = 0.1851180661871173
$ws.Range("P2").Value = 0.1851180661871173
$ws.Range("Q2").Value = 0.060828363075
$ws.Range("R2").Value = 0.5474552676750001
$ws.Range("S2").Value = 0.04266013131614595
$ws.Range("T2").Value = 0.04266013131614595

# Row 3
$ws.Range("G3").Value = 0.102275
$ws.Range("H3").Value = 0.306825
$ws.Range("I3").Value = 0.2304482333616488
$ws.Range("J3").Value = 0.2304482333616488
$ws.Range("O3").Value = 0.4463626675210189
$ws.Range("P3").Value = 0.4463626675210189
$ws.Range("Q3").Value = 0.1466713161083333
$ws.Range("R3").Value = 1.320041844975
$ws.Range("S3").Value = 0.1028634881688118
$ws.Range("T3").Value = 0.1028634881688118

# Row 4
$ws.Range("G4").Value = 0.102275
$ws.Range("H4").Value = 0.306825
$ws.Range("I4").Value = 0.2304482333616488
$ws.Range("J4").Value = 0.2304482333616488
$ws.Range("M4").Value = 1.18399
$ws.Range("N4").Value = 3.55197
$ws.Range("O4").Value = 0.3685192662918639
$ws.Range("P4").Value = 0.3685192662918639
$ws.Range("Q4").Value = 0.12109257725
$ws.Range("R4").Value = 1.08983319525
$ws.Range("S4").Value = 0.08492461387669106
$ws.Range("T4").Value = 0.08492461387669105

# Row 5
$ws.Range("I5").Value = 0.7695517666383512
$ws.Range("J5").Value = 0.7695517666383511
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.594753
$ws.Range("N5").Value = 1.784259
$ws.Range("O5").Value = 0.1851180661871173
$ws.Range("P5").Value = 0.1851180661871173
$ws.Range("Q5").Value = 0.203128371102
$ws.Range("R5").Value = 1.828155339918
$ws.Range("S5").Value = 0.1424579348709713
$ws.Range("T5").Value = 0.1424579348709713

# Row 6
$ws.Range("I6").Value = 0.7695517666383512
$ws.Range("J6").Value = 0.7695517666383511
$ws.Range("O6").Value = 0.4463626675210189
$ws.Range("P6").Value = 0.4463626675210189
$ws.Range("S6").Value = 0.3434991793522071
$ws.Range("T6").Value = 0.343499179352207

# Row 7
$ws.Range("I7").Value = 0.7695517666383512
$ws.Range("J7").Value = 0.7695517666383511
$ws.Range("M7").Value = 1.18399
$ws.Range("N7").Value = 3.55197
$ws.Range("O7").Value = 0.3685192662918639
$ws.Range("P7").Value = 0.3685192662918639
$ws.Range("Q7").Value = 0.40437284066
$ws.Range("R7").Value = 3.63935556594
$ws.Range("S7").Value = 0.2835946524151729
$ws.Range("T7").Value = 0.2835946524151728
